$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 6, shifting existing rows 6-8 down to 7-9.
$ws.Rows.Item(6).Insert()

# Fill the new row 6 with the data for the new weekly entry (copy of row 7's
# static descriptive columns, with the new date/volume/price/precio-kg values).
$ws.Range("A6").Value = 7
$ws.Range("B6").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C6").Value = "Ñuble"
$ws.Range("D6").Value = 45086
$ws.Range("D6").NumberFormat = $ws.Range("D7").NumberFormat
$ws.Range("E6").Value = 16
$ws.Range("F6").Value = "Fruta"
$ws.Range("G6").Value = 100107
$ws.Range("H6").Value = "Otros"
$ws.Range("I6").Value = 100107001
$ws.Range("J6").Value = "Caqui"
$ws.Range("K6").Value = "Mankaki"
$ws.Range("L6").Value = "Primera"
$ws.Range("M6").Value = 30
$ws.Range("N6").Value = 18000
$ws.Range("O6").Value = 18000
$ws.Range("P6").Value = 18000
$ws.Range("Q6").Value = "$/caja 18 kilos granel"
$ws.Range("R6").Value = "Región del Maule"
$ws.Range("S6").Value = 1000
$ws.Range("T6").Value = 18
